$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of Price cells carry an intentional trailing zero (e.g. "157.20",
# "0.400") that must survive as literal text -- otherwise Excel's normal
# numeric auto-detection on assignment would coerce them to numbers and the
# General format would silently drop the trailing zero (157.2, 0.4, ...).
# Marking just those cells as Text before writing keeps everything else
# untouched (no style churn on the ~84 other edited cells).
foreach ($addr in @("D6", "D11", "D34", "D46", "D50")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '65.731.72'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").Value = '2.674.35'
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '600.44'
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("D6").Value = '157.20'
$ws.Range("E6").Value = '  +0.38%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  +4.13%  '
$ws.Range("D9").Value = '0.123'
$ws.Range("E9").Value = '  -0.68%  '
$ws.Range("D10").Value = '5.93'
$ws.Range("E10").Value = '  +1.42%  '
$ws.Range("D11").Value = '0.400'
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("E12").Value = '  -0.26%  '
$ws.Range("D13").Value = '29.61'
$ws.Range("E13").Value = '  -0.35%  '
$ws.Range("D14").Value = '0.0000197'
$ws.Range("E14").Value = '  +1.90%  '
$ws.Range("E15").Value = '  +1.01%  '
$ws.Range("D16").Value = '65.504.64'
$ws.Range("D17").Value = '2.673.27'
$ws.Range("E17").Value = '  +1.54%  '
$ws.Range("D18").Value = '12.58'
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("D19").Value = '4.83'
$ws.Range("E19").Value = '  -0.65%  '
$ws.Range("E20").Value = '  +2.17%  '
$ws.Range("D21").Value = '352.01'
$ws.Range("E21").Value = '  -1.32%  '
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value = '69.75'
$ws.Range("E23").Value = '  +0.48%  '
$ws.Range("D24").Value = '0.0000111'
$ws.Range("E24").Value = '  +5.91%  '
$ws.Range("D25").Value = '9.78'
$ws.Range("E25").Value = '  +4.24%  '
$ws.Range("E26").Value = '  -3.82%  '
$ws.Range("E27").Value = '  +1.93%  '
$ws.Range("D28").Value = '1.61'
$ws.Range("E28").Value = '  -0.92%  '
$ws.Range("D29").Value = '8.17'
$ws.Range("E29").Value = '  +0.75%  '
$ws.Range("D30").Value = '542.38'
$ws.Range("E30").Value = '  +2.23%  '
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("D32").Value = '2.15'
$ws.Range("E32").Value = '  -1.15%  '
$ws.Range("D33").Value = '1.78'
$ws.Range("E33").Value = '  +0.10%  '
$ws.Range("D34").Value = '6.60'
$ws.Range("E34").Value = '  +4.31%  '
$ws.Range("D35").Value = '5.48'
$ws.Range("E35").Value = '  -1.09%  '
$ws.Range("D36").Value = '0.424'
$ws.Range("E36").Value = '  -1.57%  '
$ws.Range("E37").Value = '  -0.54%  '
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '1.96'
$ws.Range("E39").Value = '  -0.92%  '
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = '157.92'
$ws.Range("E40").Value = '  -2.32%  '
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D42").Value = '42.62'
$ws.Range("E42").Value = '  +1.42%  '
$ws.Range("D43").Value = '165.66'
$ws.Range("E43").Value = '  +0.91%  '
$ws.Range("E44").Value = '  -1.31%  '
$ws.Range("D45").Value = '0.0616'
$ws.Range("E45").Value = '  +1.44%  '
$ws.Range("D46").Value = '2.30'
$ws.Range("E46").Value = '  -2.85%  '
$ws.Range("D47").Value = '23.31'
$ws.Range("E47").Value = '  +2.27%  '
$ws.Range("D48").Value = '0.647'
$ws.Range("E48").Value = '  -0.64%  '
$ws.Range("E49").Value = '  -0.19%  '
$ws.Range("D50").Value = '0.100'
$ws.Range("E50").Value = '  +2.31%  '
$ws.Range("D51").Value = '20.05'
$ws.Range("E51").Value = '  +2.36%  '
